$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. In the source workbook every cell
# in this table is stored as text (inline string), even the ones that look
# like numbers (e.g. "243.81", "0.05601"). To keep them as text instead of
# letting Excel reinterpret them as numeric values, each cell is switched to
# the Text number format ("@") before the new value is assigned.
$changes = [ordered]@{
    "D2"  = "243.81"
    "D3"  = "21.50"
    "D4"  = "5.225"
    "D5"  = "0.05601"
    "D7"  = "6.387"
    "D8"  = "0.8054"
    "D9"  = "0.9804"
    "D10" = "0.0005746"
    "E10" = "9OneONE"
    "D11" = "0.1413"
    "D12" = "0.07304"
    "D13" = "0.03107"
    "D14" = "0.03058"
    "D15" = "0.09287"
    "D16" = "3.564"
    "D17" = "0.001644"
    "D18" = "0.04708"
    "D19" = "0.006382"
    "D20" = "0.004989"
    "D21" = "0.001042"
    "D24" = "3.759"
    "D25" = "2.098"
    "D26" = "0.3261"
    "D40" = "0.03913"
    "D41" = "0.006890"
    "D42" = "0.003398"
    "D43" = "0.1035"
    "D44" = "0.008499"
    "D45" = "0.00005930"
    "D47" = "0.0005497"
    "D48" = "0.6820"
    "D49" = "0.08919"
    "E49" = "48BOLOBOLOBestin24h"
    "D50" = "0.00002099"
}

foreach ($addr in $changes.Keys) {
    # Ensure text formatting so Excel keeps the value as a string (preserves
    # leading/trailing zeros like "21.50" or "0.0005746") rather than
    # reinterpreting it as a number.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
}
